$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped
# from 45190 to 45192 for every data row (rows 2-98).
$ws.Range("C2:C98").Value = 45192
